# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing "Late" / heading / "Outstanding" columns one
# position to the right, and make "Repayment schedule" the active sheet/tab.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N").Insert()

# Leave the selection where Excel would land it after the insert.
$wsSchedule.Range("R6").Select()

# Activate the "Repayment schedule" sheet last so it becomes the workbook's
# active tab (matching activeTab="2" / tabSelected moving off NewLoanInput).
$wsSchedule.Activate()
